$wb = $excel.ActiveWorkbook

# ----- Overall sheet -----
$wsOverall = $wb.Worksheets.Item("Overall")

$wsOverall.Range("B2").Value = 210
$wsOverall.Range("C2").Value = 7
$wsOverall.Range("D2").Value = 1.231902869808319
$wsOverall.Range("E2").Value = 0.29273040859629668
$wsOverall.Range("F2").Value = 1.5385714285714283
$wsOverall.Range("G2").Value = 342
$wsOverall.Range("H2").Value = 70
$wsOverall.Range("I2").Value = 412
$wsOverall.Range("J2").Value = 1097
$wsOverall.Range("K2").Value = 216

# ----- Zones sheet -----
$wsZones = $wb.Worksheets.Item("Zones")

# Row 2 (Zone 1)
$wsZones.Range("C2").Value = 4
$wsZones.Range("D2").Value = 0.45333333333333325
$wsZones.Range("E2").Value = 0.43981481481481466
$wsZones.Range("F2").Value = 0.57500000000000007

# Row 3 (Zone 2)
$wsZones.Range("B3").Value = 25
$wsZones.Range("D3").Value = 1.221759259259259
$wsZones.Range("E3").Value = 0.26851851851851838
$wsZones.Range("F3").Value = 1.5395061728395061

# Row 4 (Zone 3)
$wsZones.Range("B4").Value = 22
$wsZones.Range("D4").Value = 0.66896551724137943
$wsZones.Range("E4").Value = 0.2777777777777779
$wsZones.Range("F4").Value = 0.84500000000000008

# Row 5 (Zone 4)
$wsZones.Range("B5").Value = 8
$wsZones.Range("D5").Value = 1.9944444444444445
$wsZones.Range("E5").Value = 0.3
$wsZones.Range("F5").Value = 2.3333333333333335

# Row 6 (Zone 5)
$wsZones.Range("B6").Value = 23
$wsZones.Range("D6").Value = 0.3520833333333333
$wsZones.Range("E6").Value = 0.15555555555555545
$wsZones.Range("F6").Value = 0.47000000000000003

# Row 7 (Zone 6)
$wsZones.Range("B7").Value = 15
$wsZones.Range("D7").Value = 1.1746666666666667
$wsZones.Range("E7").Value = 0.11666666666666625
$wsZones.Range("F7").Value = 1.3761904761904762

# Row 8 (Zone 7)
$wsZones.Range("B8").Value = 20
$wsZones.Range("D8").Value = 1.6672222222222222
$wsZones.Range("E8").Value = 0.16111111111111118
$wsZones.Range("F8").Value = 2.0437499999999997

# Row 9 (Zone 8)
$wsZones.Range("B9").Value = 14
$wsZones.Range("D9").Value = 1.1446666666666667
$wsZones.Range("E9").Value = 0.13888888888888865
$wsZones.Range("F9").Value = 1.281818181818182

# Row 10 (Zone 9)
$wsZones.Range("B10").Value = 11
$wsZones.Range("C10").Value = 3
$wsZones.Range("D10").Value = 2.4153846153846148
$wsZones.Range("E10").Value = 0.48095238095238096
$wsZones.Range("F10").Value = 3.1280701754385958

# Row 11 (Zone 10)
$wsZones.Range("B11").Value = 4
$wsZones.Range("D11").Value = 2.5348484848484856
$wsZones.Range("E11").Value = 0.46666666666666656
$wsZones.Range("F11").Value = 2.741666666666668

# Row 12 (Zone 11)
$wsZones.Range("B12").Value = 22
$wsZones.Range("D12").Value = 0.29393939393939378
$wsZones.Range("E12").Value = 0.10555555555555533
$wsZones.Range("F12").Value = 0.32368421052631569

# Row 13 (Zone 12)
$wsZones.Range("B13").Value = 10
$wsZones.Range("D13").Value = 1.5356060606060604
$wsZones.Range("E13").Value = 0.34999999999999964
$wsZones.Range("F13").Value = 1.6541666666666668

# Row 14 (Zone 13)
$wsZones.Range("B14").Value = 20
$wsZones.Range("D14").Value = 0.64679487179487194
$wsZones.Range("E14").Value = 0.23666666666666664
$wsZones.Range("F14").Value = 0.74444444444444458
